# Clean up code and fix output
#
# Adds a new "Yearly demand" worksheet as the final tab of the workbook.
# It mirrors the layout used by every other sheet in this model (e.g.
# "DG Dispatch"): a header row (B1:Y1) numbered 0..23, a left-hand label
# column (A2:A4) numbered 0..2, and a 3x24 block of hourly demand figures
# for three representative days.

$wb = $excel.ActiveWorkbook

# Use the first sheet as the formatting template: its header/label cells
# already carry the workbook's standard bold/centered/bordered "label"
# style, so copying their format (rather than re-building it property by
# property) keeps the new sheet on the exact same style as the rest of the
# workbook.
$template = $wb.Worksheets.Item(1)

# Insert the new sheet after the last existing sheet ("Connected
# Households") so it becomes the 14th / final tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Yearly demand"

# Match the page margins used by the rest of the workbook (0.75"/0.75"/1"/1"
# with 0.5" header/footer - PageSetup margins are expressed in points).
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Header row: 0..23 across columns B..Y
for ($col = 2; $col -le 25; $col++) {
    $ws.Cells.Item(1, $col).Value = $col - 2
}

# Row labels (A2:A4) = 0, 1, 2
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(4, 1).Value = 2

# Match the bold/centered/bordered "label" style used throughout the
# workbook by copying the format from the equivalent cells on sheet 1.
$template.Range("B1:Y1").Copy()
$ws.Range("B1:Y1").PasteSpecial(-4122)
$template.Range("A2:A4").Copy()
$ws.Range("A2:A4").PasteSpecial(-4122)

# Hourly demand data for three representative days
$row2 = @(-32.5, -19.5, -13, -13, -13, 142.5, 291.5, 327, 388.5, 502, 596, 670.5, 745, 651, 576.5, 502, 320.5, 139, 32, -117, -97.5, -78, -52, -39)
$row3 = @(-32.5, -19.5, -13, 0, 0, -19.5, 0, 324, 486, 648, 729, 751.5, 583, 567, 333.5, 340, 243, 57.99999999999999, -130, 0, 0, -78, 0, -39)
$row4 = @(-32.5, -19.5, 0, 0, 0, -19.5, 0, 0, 81, 324, 567, 589.5, 648, 567, 324, 162, 81, 0, -130, 0, 0, 0, 0, -39)

for ($i = 0; $i -lt 24; $i++) {
    $ws.Cells.Item(2, $i + 2).Value = $row2[$i]
    $ws.Cells.Item(3, $i + 2).Value = $row3[$i]
    $ws.Cells.Item(4, $i + 2).Value = $row4[$i]
}
